# Apply updated crypto price/volume figures to sheet1 (matches commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several Price (column D) values look like plain numbers (e.g. "211.65").
# The source data stores them as text, so a leading apostrophe is used to force
# Excel to keep them as text instead of converting them to numeric values
# (which would also risk floating point rounding/precision drift).

$ws.Range("D2").Value = "26.720.55"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.601.61"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'211.65"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "'0.0620"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'19.67"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "1.826.35"
$ws.Range("D13").Value = "1.605.03"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "'64.92"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "26.689.15"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "'210.19"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("E21").Value = "  +2.64%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'2.29"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").Value = "'8.95"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").Value = "'144.29"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").Value = "'15.39"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D34").Value = "1.298.52"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  +11.84%  "
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "'0.827"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "'63.04"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("D45").Value = "1.739.65"
$ws.Range("D46").Value = "'90.59"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D51").Value = "'7.44"
$ws.Range("E51").Value = "  -0.03%  "
